$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 7341
$ws1.Range("F9").Value = 599
$ws1.Range("F12").Value = 4281
$ws1.Range("F13").Value = 1736
$ws1.Range("F16").Value = 2884
$ws1.Range("F20").Value = 476
$ws1.Range("F22").Value = 449
$ws1.Range("F23").Value = 287
$ws1.Range("F24").Value = 89
$ws1.Range("F25").Value = 1676
$ws1.Range("F26").Value = 1159
$ws1.Range("F27").Value = 88
$ws1.Range("F28").Value = 1360
$ws1.Range("F36").Value = 53
$ws1.Range("F37").Value = 2800
$ws1.Range("F38").Value = 696
$ws1.Range("F39").Value = 20

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 7341
$ws4.Range("F9").Value = 599
$ws4.Range("F12").Value = 4281
$ws4.Range("F13").Value = 1736
$ws4.Range("F16").Value = 2884
$ws4.Range("F20").Value = 476
$ws4.Range("F22").Value = 449
$ws4.Range("F23").Value = 287
$ws4.Range("F24").Value = 89
$ws4.Range("F25").Value = 1676
$ws4.Range("F26").Value = 1159
$ws4.Range("F27").Value = 88
$ws4.Range("F28").Value = 1360
$ws4.Range("F36").Value = 53
$ws4.Range("F37").Value = 2800
$ws4.Range("F39").Value = 696
$ws4.Range("F40").Value = 20
